$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Reorder the "Periodo Mora" column (E16:E22) from ascending to descending
# (2403, 2402, 2401, 2312, 2311, 2310, 2309) and move the corresponding
# "Valor Mora" figure (29387) so it stays attached to period 2403.
$ws.Range("E16").Value = "2403"
$ws.Range("E17").Value = "2402"
$ws.Range("E18").Value = "2401"
$ws.Range("E19").Value = "2312"
$ws.Range("E20").Value = "2311"
$ws.Range("E21").Value = "2310"
$ws.Range("E22").Value = "2309"

$ws.Range("F16").Value = 29387
$ws.Range("F22").Value = 46400
